# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) on Hoja1 listed the four arrears periods
# in descending order (1903, 1902, 1901, 1812). The update re-sorts them
# in ascending order (1812, 1901, 1902, 1903) for the four worker rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1812"
$ws.Range("E17").Value = "1901"
$ws.Range("E18").Value = "1902"
$ws.Range("E19").Value = "1903"
